$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update dimension/measure labels for sector-descripcion, aragon(->refArea), regimen, direccion-provincial-nombre
$ws.Range("E2").Value = "iaest-measure:sector-descripcion"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "iaest-measure:regimen"
$ws.Range("H2").Value = "iaest-measure:direccion-provincial-nombre"

# Row 3: sector-descripcion, regimen, direccion-provincial-nombre become "medida" (aragon/refArea column stays "dim")
$ws.Range("E3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "medida"

# Row 4: sector-descripcion, regimen, direccion-provincial-nombre become "xsd:int"; aragon/refArea column becomes "URI-Comunidad"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"

# Row 5 (mapping file references) is removed entirely
$ws.Rows.Item(5).Delete()
